$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "endTime" test-case values (rows 10/11) are entered as quote-prefixed text
# (leading apostrophe), same as typing '10/7/20 22:0x directly into Excel, so the
# stored value stays a plain date-looking string and the cell keeps a text format.
# Touch J10 first (placeholder) so its own (border-less) quote-prefix style is
# allocated before J11's (bordered) quote-prefix style, then fill in the real
# endTime values, J11 before J10, to match the shared-string build order.
$ws.Range("J10").Value = "'x"
$ws.Range("J11").Value = "'10/7/20 22:00"
$ws.Range("J10").Value = "'10/7/20 22:01"

# Row 4: eventCoordinator / eventCoordinatorError test data changed from "jyotir" to "alan"
$ws.Range("K4").Value = "alan"
$ws.Range("N4").Value = "alan is not available at that time"

# Update the active selection to N4 (was G11)
$ws.Range("N4").Select()
